$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D for the two newest quarters
$ws.Columns("D:E").Insert()

# Copy number formats from the (now-shifted) old D:E columns (F:G) into the new D:E columns
$ws.Range("F5:G102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new columns D (period ending 2018-12-31) and E (period ending 2018-09-30)
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 256400
$ws.Range("E8").Value = 203300
$ws.Range("D9").Value = 77900
$ws.Range("E9").Value = 74400
$ws.Range("D10").Value = 178500
$ws.Range("E10").Value = 128900
$ws.Range("D12").Value = 46400
$ws.Range("E12").Value = 46500
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = "NA"
$ws.Range("D15").Value = 1600
$ws.Range("E15").Value = 1600
$ws.Range("D17").Value = 240900
$ws.Range("E17").Value = 220500
$ws.Range("D18").Value = 15500
$ws.Range("E18").Value = -17200
$ws.Range("D20").Value = 2400
$ws.Range("E20").Value = 1000
$ws.Range("D21").Value = 24500
$ws.Range("E21").Value = -10000
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("D23").Value = 17900
$ws.Range("E23").Value = -16200
$ws.Range("D24").Value = 1500
$ws.Range("E24").Value = -8600
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 16400
$ws.Range("E26").Value = -7600
$ws.Range("D27").Value = 16400
$ws.Range("E27").Value = -7600
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -2400
$ws.Range("E32").Value = -1000
$ws.Range("D33").Value = 16400
$ws.Range("E33").Value = -7600
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 16400
$ws.Range("E35").Value = -7600
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 114400
$ws.Range("E41").Value = 106200
$ws.Range("D42").Value = 93000
$ws.Range("E42").Value = 99800
$ws.Range("D43").Value = 357200
$ws.Range("E43").Value = 309600
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 46000
$ws.Range("E45").Value = 70600
$ws.Range("D46").Value = 610600
$ws.Range("E46").Value = 586100
$ws.Range("D47").Value = 153800
$ws.Range("E47").Value = 168900
$ws.Range("D48").Value = "NA"
$ws.Range("E48").Value = "NA"
$ws.Range("D49").Value = 94100
$ws.Range("E49").Value = 96300
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 124000
$ws.Range("E52").Value = 111300
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 982600
$ws.Range("E54").Value = 962600
$ws.Range("D57").Value = 16500
$ws.Range("E57").Value = 12900
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("D59").Value = 315300
$ws.Range("E59").Value = 269300
$ws.Range("D60").Value = 331800
$ws.Range("E60").Value = 282300
$ws.Range("D61").Value = 0
$ws.Range("E61").Value = 0
$ws.Range("D62").Value = 29200
$ws.Range("E62").Value = 59500
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 361000
$ws.Range("E66").Value = 341800
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 510900
$ws.Range("E72").Value = 496800
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 621500
$ws.Range("E76").Value = 620900
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 16400
$ws.Range("E81").Value = -7600
$ws.Range("D83").Value = 6600
$ws.Range("E83").Value = 6200
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 37300
$ws.Range("E89").Value = -8300
$ws.Range("D91").Value = "NA"
$ws.Range("E91").Value = "NA"
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = 1400
$ws.Range("E94").Value = -3200
$ws.Range("D96").Value = -2400
$ws.Range("E96").Value = -2400
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -29800
$ws.Range("E100").Value = -25800
$ws.Range("D101").Value = -600
$ws.Range("E101").Value = -700
$ws.Range("D102").Value = 8200
$ws.Range("E102").Value = -38100

# Restate the Dec-2017 / Sep-2017 figures (now columns H and I) with corrected figures
$ws.Range("H8").Value = 254600
$ws.Range("I8").Value = 191000
$ws.Range("H10").Value = 180600
$ws.Range("I10").Value = 121300
$ws.Range("H17").Value = 216100
$ws.Range("I17").Value = 193300
$ws.Range("H18").Value = 38500
$ws.Range("I18").Value = -2300
$ws.Range("H20").Value = -1200
$ws.Range("I20").Value = -5000
$ws.Range("H21").Value = 43300
$ws.Range("I21").Value = -900
$ws.Range("H23").Value = 37300
$ws.Range("I23").Value = -7200
$ws.Range("H24").Value = 10900
$ws.Range("I24").Value = -8500
$ws.Range("H26").Value = 26400
$ws.Range("I26").Value = 1300
$ws.Range("H27").Value = 26400
$ws.Range("I27").Value = 1300
$ws.Range("H29").Value = 14200
$ws.Range("H32").Value = 1200
$ws.Range("I32").Value = 5000
$ws.Range("H33").Value = 40600
$ws.Range("I33").Value = 1300
$ws.Range("H35").Value = 40600
$ws.Range("I35").Value = 1300
$ws.Range("H43").Value = 657700
$ws.Range("H46").Value = 647700
$ws.Range("H52").Value = 163600
$ws.Range("H54").Value = 1012800
$ws.Range("H59").Value = 277400
$ws.Range("H60").Value = 294800
$ws.Range("H62").Value = 88700
$ws.Range("H66").Value = 356900
$ws.Range("H72").Value = 509700
$ws.Range("H76").Value = 655900
$ws.Range("H81").Value = 40600
$ws.Range("I81").Value = 1300
